$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 (H) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 298
$wsOff.Range("C2").Value = 201
$wsOff.Range("D2").Value = 69
$wsOff.Range("E2").Value = 24

# --- DEF sheet: row 2 (H) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 263
$wsDef.Range("C2").Value = 159
$wsDef.Range("D2").Value = 45
$wsDef.Range("E2").Value = 14
$wsDef.Range("F2").Value = 4
$wsDef.Range("G2").Value = 4
